# Update "want to go" counts (F column) across the relevant sheets.
# These values are regenerated by the site build (bilibili ticket counts)
# at each gh-pages deploy, so only the numeric values change.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 6281
$wsExhibition.Range("F3").Value = 566
$wsExhibition.Range("F4").Value = 123
$wsExhibition.Range("F8").Value = 1383

# Sheet "演出" (Performance)
$wsPerformance = $wb.Worksheets.Item("演出")
$wsPerformance.Range("F2").Value = 98

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 6281
$wsAll.Range("F3").Value = 566
$wsAll.Range("F4").Value = 123
$wsAll.Range("F8").Value = 98
$wsAll.Range("F12").Value = 1383
